# Append new evaluation-run results to each of the three result sheets,
# mirroring additional runs of the (now improved) classification algorithms.

$wb = $excel.ActiveWorkbook

# --- Sheet "NaiveBayes": add rows 8-12 -------------------------------------
$ws1 = $wb.Worksheets.Item("NaiveBayes")
$naiveBayesRows = @(
    @("09/11/2022 09:28:43", 0.545),
    @("09/11/2022 09:32:53", 0.545),
    @("09/11/2022 09:33:17", 0.545),
    @("09/11/2022 09:35:40", 0.545),
    @("09/11/2022 09:37:12", 0.545)
)
$r = 8
foreach ($row in $naiveBayesRows) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# --- Sheet "LogisticRegression": add rows 7-10 ------------------------------
$ws2 = $wb.Worksheets.Item("LogisticRegression")
$logRegRows = @(
    @("09/11/2022 09:28:46", 0.475),
    @("09/11/2022 09:33:20", 0.58125),
    @("09/11/2022 09:35:43", 0.585),
    @("09/11/2022 09:37:15", 0.585)
)
$r = 7
foreach ($row in $logRegRows) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# --- Sheet "SVM": add rows 4-6 ----------------------------------------------
$ws3 = $wb.Worksheets.Item("SVM")
$svmRows = @(
    @("09/11/2022 09:33:39", 0.4825),
    @("09/11/2022 09:36:03", 0.4825),
    @("09/11/2022 09:37:35", 0.59)
)
$r = 4
foreach ($row in $svmRows) {
    $ws3.Cells.Item($r, 1).Value = $row[0]
    $ws3.Cells.Item($r, 2).Value = $row[1]
    $r++
}
